$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove every paragraph between the chapter title (paragraph 1)
#    and the final (empty) paragraph that precedes the sectPr. This
#    collapses the whole old "CONCLUSIONES Y RECOMENDACIONES" body
#    down to just those two anchor paragraphs.
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
if ($count -gt 2) {
    $start = $d.Paragraphs(2).Range.Start
    $end = $d.Paragraphs($count - 1).Range.End
    $midRange = $d.Range($start, $end)
    $midRange.Delete()
}

# ------------------------------------------------------------------
# 2. Rewrite the chapter title paragraph: "CAPITULO XI - CONCLUSIONES
#    Y RECOMENDACIONES" -> "CAPITULO IX - CUARTA ITERACÍON"
# ------------------------------------------------------------------
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo1"/></w:pPr><w:bookmarkStart w:id="0" w:name="_Toc431546804"/><w:r><w:t xml:space="preserve">CAPITULO </w:t></w:r><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>IX</w:t></w:r><w:r><w:t xml:space="preserve"> – </w:t></w:r><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>CUARTA</w:t></w:r><w:r><w:t xml:space="preserve"> ITERACÍON</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p>'
$d.Paragraphs(1).Range.InsertXML($titleXml)

# ------------------------------------------------------------------
# 3. Replace the trailing empty paragraph with the new chapter body:
#    an intro paragraph plus the 6.x sub-heading skeleton, followed
#    by one bare empty paragraph before the section break.
# ------------------------------------------------------------------
$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/><w:rPr><w:b w:val="0"/><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">La metodologia scrum es iterativa e incremental, por lo tanto en este capitulo, se abarca todo lo relacionado al desarrollo de la </w:t></w:r><w:r><w:rPr><w:b w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t>cuarta</w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/><w:r><w:rPr><w:b w:val="0"/><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> iteracion.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t>6.1 Planificación de la iteración</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t>6.2 Pila de la Iteración</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t>6.3 Diseño</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:t>6.5 Pruebas</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>6.6 Grafico BurnDown</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Ttulo2"/></w:pPr><w:r><w:t>6.</w:t></w:r><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>7</w:t></w:r><w:r><w:t xml:space="preserve"> Revisión de la Iteración</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>'
$d.Paragraphs($d.Paragraphs.Count).Range.InsertXML($bodyXml)

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
